# refactor: CheckpointManager and add checkpointing to get_candidate_data
#
# Replaces the old "medication renaming" sample rows (items 11,13,14,16,17,18)
# with new "SORO" (saline solution) rows (items 6,7,8,9), shrinking the table
# from A1:H7 down to A1:H5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-5 (row 1 header stays the same)
$data = @(
    @(6, "SORO GLICOFISIOLÓGICO, GLICOSE À 5% EM CLORETO DE SÓDIO À 0,9%; BOLSA SISTEMA FECHADO, FRASCO COM 1000ML", "CLORETO DE SÓDIO;GLICOSE", "1000ml", "JP", "JP INDUSTRIA FARMACEUTICA S/A"),
    @(7, "SORO GLICOFISIOLÓGICO, GLICOSE À 5% EM CLORETO DE SÓDIO À 0,9%; BOLSA SISTEMA FECHADO, FRASCO COM 250ML", "CLORETO DE SÓDIO;GLICOSE", "250ml", "JP", "JP INDUSTRIA FARMACEUTICA S/A"),
    @(8, "SORO GLICOFISIOLÓGICO, GLICOSE À 5% EM CLORETO DE SÓDIO À 0,9%; BOLSA SISTEMA FECHADO, FRASCO COM 500ML", "CLORETO DE SÓDIO;GLICOSE", "500ml", "JP", "JP INDUSTRIA FARMACEUTICA S/A"),
    @(9, "SORO GLICOSADO 5%, BOLSA SISTEMA FECHADO, FRASCO COM 1000ML", "SORO GLICOSADO 5%, BOLSA SISTEMA FECHADO, FRASCO COM 1000ML", "1000ml", "JP", "JP INDUSTRIA FARMACEUTICA S/A")
)

# Remove the two extra rows (old table had 6 data rows, new table has 4)
$ws.Range("A7:H7").EntireRow.Delete()
$ws.Range("A6:H6").EntireRow.Delete()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
}
